$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column D (rows 2-25), each equal to old value + 16.8
$values = @(16.4, 13.85, 11.31, 8.76, 6.22, 3.67, 16.4, 13.85, 11.31, 8.76, 6.22, 3.67, 16.4, 13.85, 11.31, 8.76, 6.22, 3.67, 16.4, 13.85, 11.31, 8.76, 6.22, 3.67)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $values[$i]
}
